# Adds a new "2022-Q3" quarter sheet (copied/restyled from "2022-Q2") with
# its own fund-holding data, and records the new quarter in the "总计"
# (totals) summary sheet. All other quarter sheets are left untouched and
# simply shift right in tab order because the new sheet is inserted before
# them.

function Set-CellText($cell, [string]$text) {
    # Cells that "look numeric" (digits/dot only, e.g. fund codes like
    # "012744" or figures like "5.60") get silently coerced to a Number by
    # plain `.Value =` assignment (losing leading zeros / trailing zeros).
    # Forcing a text number-format first keeps them as literal text, same
    # as the inlineStr cells in the source workbook.
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it keeps
#    identical column layout / header styling / border formatting), place
#    it immediately before "2022-Q2", then overwrite it with the new
#    quarter's fund data.
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($sheetQ2)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Row 2: 012744 / 光大保德信品质生活混合A
Set-CellText $newSheet.Cells.Item(2,2) "012744"
Set-CellText $newSheet.Cells.Item(2,3) "光大保德信品质生活混合A"
Set-CellText $newSheet.Cells.Item(2,4) "5.60"
Set-CellText $newSheet.Cells.Item(2,5) "88.62"
Set-CellText $newSheet.Cells.Item(2,6) "5.58"
Set-CellText $newSheet.Cells.Item(2,7) "0.3125"
$newSheet.Cells.Item(2,8).Value = 5

# Row 3: 007592 / 华夏价值精选混合 (name/code unchanged from 2022-Q2)
Set-CellText $newSheet.Cells.Item(3,2) "007592"
Set-CellText $newSheet.Cells.Item(3,3) "华夏价值精选混合"
Set-CellText $newSheet.Cells.Item(3,4) "2.26"
Set-CellText $newSheet.Cells.Item(3,5) "93.77"
Set-CellText $newSheet.Cells.Item(3,6) "6.34"
Set-CellText $newSheet.Cells.Item(3,7) "0.1433"
$newSheet.Cells.Item(3,8).Value = 5

# Row 4: 012758 / 光大保德信品质生活混合C
Set-CellText $newSheet.Cells.Item(4,2) "012758"
Set-CellText $newSheet.Cells.Item(4,3) "光大保德信品质生活混合C"
Set-CellText $newSheet.Cells.Item(4,4) "0.35"
Set-CellText $newSheet.Cells.Item(4,5) "88.62"
Set-CellText $newSheet.Cells.Item(4,6) "5.58"
Set-CellText $newSheet.Cells.Item(4,7) "0.0195"
$newSheet.Cells.Item(4,8).Value = 5

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 row at the top of
#    the data (row 2) and push the existing quarters down by one row. The
#    index column (A) is left exactly as-is (0,1,2,3,...) since it is
#    simply the row's position in the list.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$oldData = @()
for ($r = 2; $r -le 7; $r++) {
    $b = $total.Cells.Item($r,2).Value()
    $c = $total.Cells.Item($r,3).Value()
    $d = $total.Cells.Item($r,4).Value()
    $oldData += ,@($b, $c, $d)
}

$newRow = ,@("2022-Q3", 3, 0.48)
$allRows = $newRow + $oldData
$lastRow = $allRows.Count + 1

# The bottom row (8) is brand new - copy the bold/bordered index-column
# style from row 2 onto it before the loop below overwrites its value.
$total.Cells.Item(2, 1).Copy()
$total.Cells.Item($lastRow, 1).PasteSpecial(-4122)

for ($i = 0; $i -lt $allRows.Count; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $allRows[$i][0]
    $total.Cells.Item($r, 3).Value = $allRows[$i][1]
    $total.Cells.Item($r, 4).Value = $allRows[$i][2]
}

Write-Output "2022-Q3 sheet added and 总计 updated"
